$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.521.96"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.875.51"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.95"
$ws.Range("E5").Value = "  -3.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4875"
$ws.Range("E7").Value = "  -1.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2897"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06662"
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("D10").Value = "1.874.35"
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.52"
$ws.Range("E11").Value = "  -2.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07235"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "88.57"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.989"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6485"
$ws.Range("E15").Value = "  -3.40%  "
$ws.Range("D16").Value = "30.471.23"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007844"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.96"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").Value = "2.116.04"
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.0000"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.709"
$ws.Range("E22").Value = "  -2.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "195.48"
$ws.Range("E23").Value = "  +11.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.107"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.357"
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.72"
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.45"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.826"
$ws.Range("E28").Value = "  -5.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.408"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.247"
$ws.Range("E30").Value = "  -2.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09022"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.913"
$ws.Range("E32").Value = "  -2.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05109"
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7200"
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("E35").Value = "  -4.77%  "
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01812"
$ws.Range("E37").Value = "  -3.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.661"
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9169"
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.039"
$ws.Range("E40").Value = "  -6.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4389"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.74"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9949"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.710"
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1326"
$ws.Range("E45").Value = "  -2.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.360"
$ws.Range("E46").Value = "  -3.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4019"
$ws.Range("E47").Value = "  +3.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05823"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.611"
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.399"
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.07"
$ws.Range("E51").Value = "  -1.02%  "
